# Applies the edits described by the commit diff to the active document.
$d = $word.ActiveDocument

# 1) Sender name
$d.Content.Find.Execute("Clivaz Loris", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Burgener Sascha", 2)

# 2) Sender street address (first occurrence)
$d.Content.Find.Execute("chemin de turtemean 20", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Bachstrasse 29", 2)

# 3) Sender postal code / city
$d.Content.Find.Execute("3973 venthone", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "3931", 2)

# 4) Recipient company name
$d.Content.Find.Execute("Sierre-Energie ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "1997 ", 2)

# 5) Recipient street address
$d.Content.Find.Execute("chemin de turtemean 4", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Bachstrasse 29", 2)

# 6) Recipient postal code / city
$d.Content.Find.Execute("3973venthone  ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "3931 ", 2)

# 7) Dateline location
$d.Content.Find.Execute("venthone, ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "dfd, ", 2)

# 8) Subject line first word
$d.Content.Find.Execute("Surbooking", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Bagages", 2)

# 9) Flight number references ("256" -> "dfdf", appears twice, replace all)
$d.Content.Find.Execute("256", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "dfdf", 2)

# 10) Purchase date
$d.Content.Find.Execute("12.01.2020", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "dfd", 2)

# 11) Destination
$d.Content.Find.Execute("malte", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Lalden", 2)

# 12) Fill in the previously empty explanatory paragraph (paraId 15C8A2ED,
#     the 19th paragraph in the document).
$p = $d.Paragraphs.Item(19)
$p.Range.Text = "A mon arrivée à destination, j’ai eu la désagréable surprise de constater que la valise contenant mes effets personnels était endommagée. `nComme j’ai signalé dans les délais la chose au bureau compétent (cf. attestation en annexe)  j’ai droit à l’indemnisation du préjudice financier subi conformément à la Convention de Montréal (plafond maximum environ 1500 fr). Vous trouverez tous les justificatifs et pièces utiles en annexe.`n"

# 13) Requested indemnity amount
$d.Content.Find.Execute("2000", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "3", 2)

# 14) Bank account number
$d.Content.Find.Execute("253514564984456.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "32434.", 2)

Write-Host "Edits applied"
